$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BR6's summary text (row 7, column D): "no validation on email " -> "no validation on email on registration"
$ws.Range("D7").Value = "no validation on email on registration"

# Add "Link to attachments" hyperlink for BR5 (row 6) in column J, matching the
# style used by the other "Link to attachments" cells (J3:J5).
$ws.Hyperlinks.Add($ws.Range("J6"), "https://github.com/Oleksandr-Mnk/Test-documentation/tree/main/Bug%20reports/Attachments%20to%20bug%20reports/BR5")
$ws.Range("J6").Value = "Link to attachments"
$ws.Range("J3").Copy()
$ws.Range("J6").PasteSpecial(-4122)

# Update the view: scroll back to the top and move the selection to D8.
$ws.Application.GoTo($ws.Range("A1"))
[void]$ws.Range("D8").Select()

Write-Output "done"
